$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" -> "_FV2404", "_new" -> "_FV2410" ---
$oldCols = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldCols.Count; $i++) {
    # columns A..J -> "<name>_FV2404"
    $colLetter = [char]([int][char]'A' + $i)
    $ws.Range("$colLetter`1").Value = "$($oldCols[$i])_FV2404"

    # columns L..U -> "<name>_FV2410"
    $colLetter2 = [char]([int][char]'L' + $i)
    $ws.Range("$colLetter2`1").Value = "$($oldCols[$i])_FV2410"
}
# column K ("diff") is unchanged

# --- 2. Turn the data range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U78")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
